$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "34.098.88"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "1.789.79"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "222.59"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "32.28"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").Value = "0.0716"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "2.046.56"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.804.16"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "10.93"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "0.627"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "34.087.81"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "68.09"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "244.12"
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "10.74"
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -4.37%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "158.88"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "16.36"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  -2.99%  "
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("D34").Value = "1.80"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("D35").Value = "1.395.63"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "0.649"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "79.66"
$ws.Range("E39").Value = "  -7.11%  "
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "0.920"
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").Value = "2.16"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "0.0498"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  -2.88%  "
$ws.Range("D46").Value = "107.46"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "1.946.59"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "11.99"
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("E51").Value = "  +0.75%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
